# Rows 12-14 got rotated: the record that used to live in row 13 now lives
# in row 12, the record that used to live in row 14 now lives in row 13,
# and the record that used to live in row 12 now lives in row 14 (a
# cyclic shift). Only columns A, B, E, F, G, H, P, Q, R, S carry values
# that differ between the three source rows (the rest are identical across
# the three rows, so they don't visibly change even though the rows
# rotated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move, for rows 12-14.
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

$before = @{}
foreach ($r in 12..14) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $before[$r] = $rowVals
}

# New row 12 <- old row 13, new row 13 <- old row 14, new row 14 <- old row 12
$mapping = @{ 12 = 13; 13 = 14; 14 = 12 }

foreach ($destRow in 12..14) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $before[$srcRow][$col]
    }
}
